# Inserts a new weekly price record for "Femacal de La Calera - Espinaca"
# as row 269 of Sheet1, pushing the existing rows 269:287 down to 270:288.
# This matches the "Fruta / hortaliza, semanal" update that adds the
# latest week's observation at the top of this sub-series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 269; existing rows 269-287 shift down to 270-288.
$ws.Rows("269:269").Insert()

# Populate the new row 269 with the latest observation.
$ws.Range("A269").Value = 3
$ws.Range("B269").Value = "Femacal de La Calera"
$ws.Range("C269").Value = "Coquimbo"
$ws.Range("D269").Value = 44610
$ws.Range("E269").Value = 5
$ws.Range("F269").Value = 100112012
$ws.Range("G269").Value = "Espinaca"
$ws.Range("H269").Value = "Sin especificar"
$ws.Range("I269").Value = "Primera"
$ws.Range("J269").Value = 105
$ws.Range("K269").Value = 5000
$ws.Range("L269").Value = 5500
$ws.Range("M269").Value = 5262
$ws.Range("N269").Value = "$/docena de atados (3 kilos)"
$ws.Range("O269").Value = "Provincia de Quillota"
$ws.Range("P269").Value = 1754
$ws.Range("Q269").Value = 3
$ws.Range("R269").Value = "Hortaliza"
